# Applies the "merge import" edit: reorders the family-history columns
# J:Q on the single worksheet, moving the "famhhist_cardiomyopathy" column
# from J to L (and shifting / reordering the rest of the family-history
# columns accordingly), while keeping each column's Yes/No data attached
# to its own header.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header order for columns J..Q (text values, taken from the shared
# strings / worksheet target state).
$headers = @("famhist_none", "famhist_deaf", "famhhist_cardiomyopathy", "famhist_encephalopathy", "famhist_diabmell", "famhist_cardiovasc", "famhist_malignancy", "famhist_unknown")

# Per-header Yes/No values for data rows 2..6 (row index 0 = row2 ... 4 = row6)
$data = @{
    "famhist_none"            = @("No",  "No",  "Yes", "No",  "No")
    "famhist_deaf"             = @("No",  "Yes", "No",  "No",  "No")
    "famhhist_cardiomyopathy" = @("Yes", "Yes", "No",  "No",  "No")
    "famhist_encephalopathy"   = @("Yes", "No",  "No",  "No",  "No")
    "famhist_diabmell"         = @("Yes", "No",  "No",  "No",  "No")
    "famhist_cardiovasc"       = @("No",  "No",  "No",  "No",  "No")
    "famhist_malignancy"       = @("No",  "No",  "No",  "No",  "No")
    "famhist_unknown"          = @("No",  "No",  "No",  "No",  "Yes")
}
# fix up famhist_cardiovasc / famhist_malignancy (row5 = Yes for both in target)
$data["famhist_cardiovasc"] = @("No", "No", "No", "Yes", "No")
$data["famhist_malignancy"] = @("No", "No", "No", "Yes", "No")

$startCol = 10  # column J

for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = $startCol + $i
    $header = $headers[$i]
    $ws.Cells.Item(1, $col).Value = $header
    $vals = $data[$header]
    for ($r = 0; $r -lt $vals.Length; $r++) {
        $ws.Cells.Item(2 + $r, $col).Value = $vals[$r]
    }
}

# The "famhhist_cardiomyopathy" column carries a text-number-format style
# (numFmtId 49) that travels with it to its new home, column L (12).
$ws.Cells.Item(2, 12).NumberFormat = "@"

# Column width: the bestFit width=22 custom column width moves from column
# J (10) to column L (12) along with the relocated header.
$ws.Columns.Item(10).ColumnWidth = 18.28515625
$ws.Columns.Item(12).ColumnWidth = 22

# Update the sheet view to match the scrolled/selected state after the edit.
$window = $excel.ActiveWindow
$window.ScrollColumn = 10
$ws.Range("R15").Select()
